# Apply the "add land to growth rate" edit:
#  - update the "Last update" timestamp on the info sheet
#  - on the gases sheet, split the old CO2 row into CO2_FFI and CO2_LULUCF,
#    reorder the gas rows and add a GHG total row

$wb = $excel.ActiveWorkbook

# --- 1. Update "Last update" on the info sheet ---
$infoWs = $wb.Worksheets.Item("info")
$infoWs.Range("B2").Value = "2021-08-25 11:29:46"

# --- 2. Rewrite the gases sheet data rows (2-7) ---
$gasesWs = $wb.Worksheets.Item("gases")

# Column layout: A=category, B=subcategory (blank), C=var,
# D=1990, E=2000, F=2010, G=2019, H=1990-2000, I=2000-2010, J=2010-2019, K=1990-2019

# Row 2: CO2_FFI (previously the CO2 row's data, renamed)
$gasesWs.Cells.Item(2, 1).Value = "Gases"
$gasesWs.Cells.Item(2, 3).Value = "CO2_FFI"
$gasesWs.Cells.Item(2, 4).Value = 22.7048976097961
$gasesWs.Cells.Item(2, 5).Value = 25.8123841524056
$gasesWs.Cells.Item(2, 6).Value = 34.1369545479154
$gasesWs.Cells.Item(2, 7).Value = 37.9300547752771
$gasesWs.Cells.Item(2, 8).Value = 1.29099964213424
$gasesWs.Cells.Item(2, 9).Value = 2.83469525693971
$gasesWs.Cells.Item(2, 10).Value = 1.1775827864235
$gasesWs.Cells.Item(2, 11).Value = 1.78527699058815

# Row 3: CO2_LULUCF (new row)
$gasesWs.Cells.Item(3, 1).Value = "Gases"
$gasesWs.Cells.Item(3, 3).Value = "CO2_LULUCF"
$gasesWs.Cells.Item(3, 4).Value = 4.979180807728
$gasesWs.Cells.Item(3, 5).Value = 5.05170398532267
$gasesWs.Cells.Item(3, 6).Value = 5.33731686335467
$gasesWs.Cells.Item(3, 7).Value = 6.6050388476
$gasesWs.Cells.Item(3, 8).Value = 0.144706880948475
$gasesWs.Cells.Item(3, 9).Value = 0.551489710927444
$gasesWs.Cells.Item(3, 10).Value = 2.3961431572314
$gasesWs.Cells.Item(3, 11).Value = 0.979132921884363

# Row 4: CH4 (previously row 2's data)
$gasesWs.Cells.Item(4, 1).Value = "Gases"
$gasesWs.Cells.Item(4, 3).Value = "CH4"
$gasesWs.Cells.Item(4, 4).Value = 8.17558268994896
$gasesWs.Cells.Item(4, 5).Value = 8.43768841872268
$gasesWs.Cells.Item(4, 6).Value = 9.66424928526536
$gasesWs.Cells.Item(4, 7).Value = 10.565894553494
$gasesWs.Cells.Item(4, 8).Value = 0.316062384764004
$gasesWs.Cells.Item(4, 9).Value = 1.36650294724756
$gasesWs.Cells.Item(4, 10).Value = 0.996015124879857
$gasesWs.Cells.Item(4, 11).Value = 0.888333947826836

# Row 5: N2O (previously row 6's data)
$gasesWs.Cells.Item(5, 1).Value = "Gases"
$gasesWs.Cells.Item(5, 3).Value = "N2O"
$gasesWs.Cells.Item(5, 4).Value = 1.89826752686084
$gasesWs.Cells.Item(5, 5).Value = 2.03605111630764
$gasesWs.Cells.Item(5, 6).Value = 2.27088794788791
$gasesWs.Cells.Item(5, 7).Value = 2.53267773751827
$gasesWs.Cells.Item(5, 8).Value = 0.703166309668291
$gasesWs.Cells.Item(5, 9).Value = 1.09756671889707
$gasesWs.Cells.Item(5, 10).Value = 1.21966932213784
$gasesWs.Cells.Item(5, 11).Value = 0.999219517967798

# Row 6: Fgas (previously row 4's data)
$gasesWs.Cells.Item(6, 1).Value = "Gases"
$gasesWs.Cells.Item(6, 3).Value = "Fgas"
$gasesWs.Cells.Item(6, 4).Value = 0.286439531994979
$gasesWs.Cells.Item(6, 5).Value = 0.524520427692494
$gasesWs.Cells.Item(6, 6).Value = 0.657571963830359
$gasesWs.Cells.Item(6, 7).Value = 0.692236336011965
$gasesWs.Cells.Item(6, 8).Value = 6.23630206658874
$gasesWs.Cells.Item(6, 9).Value = 2.28644578995516
$gasesWs.Cells.Item(6, 10).Value = 0.572445739829575
$gasesWs.Cells.Item(6, 11).Value = 3.08952351408778

# Row 7: GHG (new row, totals across all gases)
$gasesWs.Cells.Item(7, 1).Value = "Gases"
$gasesWs.Cells.Item(7, 3).Value = "GHG"
$gasesWs.Cells.Item(7, 4).Value = 38.0443681663289
$gasesWs.Cells.Item(7, 5).Value = 41.8623481004511
$gasesWs.Cells.Item(7, 6).Value = 52.0669806082537
$gasesWs.Cells.Item(7, 7).Value = 58.3259022499013
$gasesWs.Cells.Item(7, 8).Value = 0.960924995466983
$gasesWs.Cells.Item(7, 9).Value = 2.20540908767668
$gasesWs.Cells.Item(7, 10).Value = 1.26926890757961
$gasesWs.Cells.Item(7, 11).Value = 1.48433332923716

Write-Host "Edit applied"
